$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.935.65"
$ws.Range("E2").Value = "  -2.68%  "

# Row 3
$ws.Range("D3").Value = "1.886.99"
$ws.Range("E3").Value = "  -3.43%  "

# Row 4
$ws.Range("E4").Value = "  -1.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.94"
$ws.Range("E5").Value = "  +1.71%  "

# Row 6
$ws.Range("E6").Value = "  -1.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4583"
$ws.Range("E7").Value = "  -3.82%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3922"
$ws.Range("E8").Value = "  -2.58%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.68"
$ws.Range("E9").Value = "  -9.84%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08223"
$ws.Range("E10").Value = "  -3.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.036"
$ws.Range("E11").Value = "  -2.32%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.76"
$ws.Range("E12").Value = "  -1.86%  "

# Row 13
$ws.Range("D13").Value = "1.892.71"
$ws.Range("E13").Value = "  -4.73%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.305"
$ws.Range("E14").Value = "  -3.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.965"
$ws.Range("E15").Value = "  -3.93%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  -1.11%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.93"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001055"
$ws.Range("E18").Value = "  -1.81%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06590"
$ws.Range("E19").Value = "  -0.31%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.46"
$ws.Range("E20").Value = "  -6.37%  "

# Row 21
$ws.Range("E21").Value = "  -1.21%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.626"
$ws.Range("E22").Value = "  -2.98%  "

# Row 23
$ws.Range("D23").Value = "27.960.39"
$ws.Range("E23").Value = "  -2.64%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  -3.92%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.300"
$ws.Range("E25").Value = "  +0.41%  "

# Row 26
$ws.Range("D26").Value = "2.137.46"
$ws.Range("E26").Value = "  -3.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.16"
$ws.Range("E27").Value = "  -0.39%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.87"
$ws.Range("E28").Value = "  -1.49%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.698"
$ws.Range("E29").Value = "  -4.23%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.102"
$ws.Range("E30").Value = "  -2.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.28"
$ws.Range("E31").Value = "  -0.40%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09525"
$ws.Range("E32").Value = "  -0.46%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9555"
$ws.Range("E33").Value = "  -4.66%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.472"
$ws.Range("E34").Value = "  +2.62%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.631"
$ws.Range("E35").Value = "  -1.09%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.453"
$ws.Range("E36").Value = "  -3.67%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.253"
$ws.Range("E37").Value = "  -0.93%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02278"
$ws.Range("E38").Value = "  -2.81%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06089"
$ws.Range("E39").Value = "  -2.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.553"
$ws.Range("E40").Value = "  -2.19%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6092"
$ws.Range("E41").Value = "  -1.86%  "

# Row 42
$ws.Range("E42").Value = "  -1.12%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.70"
$ws.Range("E43").Value = "  -3.37%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1888"
$ws.Range("E44").Value = "  -1.38%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.303"
$ws.Range("E45").Value = "  -2.08%  "

# Row 46 - Decentraland
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5807"
$ws.Range("E46").Value = "  -1.97%  "

# Row 47 - EnergySwap
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.68"
$ws.Range("E47").Value = "  -1.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.986"
$ws.Range("E48").Value = "  -4.29%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.423"
$ws.Range("E49").Value = "  +0.43%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06902"
$ws.Range("E50").Value = "  +1.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.32"
$ws.Range("E51").Value = "  -0.68%  "
